$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98, shifting existing rows 98-142 down to 99-143
$ws.Rows.Item(98).Insert()

# Populate the new row 98 with the new record's data
$ws.Range("A98").Value = 11
$ws.Range("B98").Value = "Vega Monumental Concepción"
$ws.Range("C98").Value = "Bíobío"
$ws.Range("D98").Value = 45205
$ws.Range("E98").Value = 8
$ws.Range("F98").Value = 100112012
$ws.Range("G98").Value = "Espinaca"
$ws.Range("H98").Value = "Sin especificar"
$ws.Range("I98").Value = "Primera"
$ws.Range("J98").Value = 50
$ws.Range("K98").Value = 8000
$ws.Range("L98").Value = 8000
$ws.Range("M98").Value = 8000
$ws.Range("N98").Value = "$/cuna 10 kilos"
$ws.Range("O98").Value = "Región Metropolitana"
$ws.Range("P98").Value = 800
$ws.Range("Q98").Value = 10
$ws.Range("R98").Value = "Hortaliza"
